$d = $word.ActiveDocument
$whole = $d.Content
$frag = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
  <w:pPr>
    <w:rPr>
      <w:rStyle w:val="a3"/>
      <w:rFonts w:ascii="Open Sans" w:hAnsi="Open Sans" w:cs="Open Sans"/>
      <w:color w:val="313131"/>
      <w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rStyle w:val="a3"/>
      <w:rFonts w:ascii="Open Sans" w:hAnsi="Open Sans" w:cs="Open Sans"/>
      <w:color w:val="313131"/>
      <w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>
    </w:rPr>
    <w:t>TextFile00</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:rPr>
      <w:rStyle w:val="a3"/>
      <w:rFonts w:ascii="Open Sans" w:hAnsi="Open Sans" w:cs="Open Sans"/>
      <w:color w:val="313131"/>
      <w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>
      <w:lang w:val="en-US"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rStyle w:val="a3"/>
      <w:rFonts w:ascii="Open Sans" w:hAnsi="Open Sans" w:cs="Open Sans"/>
      <w:color w:val="313131"/>
      <w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>
    </w:rPr>
    <w:t>TextFile00</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rStyle w:val="a3"/>
      <w:rFonts w:ascii="Open Sans" w:hAnsi="Open Sans" w:cs="Open Sans"/>
      <w:color w:val="313131"/>
      <w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>
      <w:lang w:val="en-US"/>
    </w:rPr>
    <w:t xml:space="preserve"> - 1</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:rPr>
      <w:lang w:val="en-US"/>
    </w:rPr>
  </w:pPr>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@
$whole.InsertXML($frag)
Write-Output ("Text=[" + $d.Content.Text + "]")
Write-Output ("ParaCount=" + $d.Paragraphs.Count)
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
  $p = $d.Paragraphs($i)
  Write-Output ($i.ToString() + ": [" + $p.Range.Text + "] Start=" + $p.Range.Start + " End=" + $p.Range.End)
}

# Fixup: InsertXML drops w:rStyle on run-level rPr (but keeps it on pPr/rPr).
# Re-apply character style "Emphasis" (a3) to each run's text range.
$p1 = $d.Paragraphs(1)
$run1Range = $d.Range($p1.Range.Start, $p1.Range.End - 1)
Write-Output ("run1Range=[" + $run1Range.Text + "]")
$run1Range.Style = "Emphasis"

$p2 = $d.Paragraphs(2)
$p2Start = $p2.Range.Start
$p2End = $p2.Range.End - 1   # exclude paragraph mark
$splitPoint = $p2Start + 10  # "TextFile00" is 10 chars
$sub1 = $d.Range($p2Start, $splitPoint)
$sub2 = $d.Range($splitPoint, $p2End)
Write-Output ("sub1=[" + $sub1.Text + "] sub2=[" + $sub2.Text + "]")
$sub1.Style = "Emphasis"
$sub2.Style = "Emphasis"

Write-Output ("Final Text=[" + $d.Content.Text + "]")
